$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) holds values that look numeric (e.g. "1.190",
# "0.00001077", "28.137.55") but must stay plain TEXT, exactly as scraped -
# matching the source inlineStr cells. Force text format before writing so
# COM does not coerce/renormalize them into real numbers.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

# --- Row-level value updates (Price / Volume columns) ---
$ws.Range("D2").Value = "28.137.55"
$ws.Range("E2").Value = "  +3.11%  "
$ws.Range("D3").Value = "1.776.67"
$ws.Range("E3").Value = "  -0.48%  "
$ws.Range("D4").Value = "1.008"
$ws.Range("E4").Value = "  +0.46%  "
$ws.Range("D5").Value = "339.33"
$ws.Range("E5").Value = "  -0.36%  "
$ws.Range("D6").Value = "1.004"
$ws.Range("E6").Value = "  +0.26%  "
$ws.Range("D7").Value = "0.3828"
$ws.Range("E7").Value = "  -4.45%  "
$ws.Range("D8").Value = "0.3428"
$ws.Range("E8").Value = "  -0.54%  "
$ws.Range("D9").Value = "46.96"
$ws.Range("E9").Value = "  -1.91%  "
$ws.Range("D10").Value = "1.152"
$ws.Range("E10").Value = "  -3.46%  "
$ws.Range("D11").Value = "0.07392"
$ws.Range("E11").Value = "  -0.82%  "
$ws.Range("D12").Value = "23.51"
$ws.Range("E12").Value = "  +8.25%  "
$ws.Range("D13").Value = "1.004"
$ws.Range("E13").Value = "  +0.33%  "
$ws.Range("D14").Value = "6.444"
$ws.Range("E14").Value = "  -0.24%  "
$ws.Range("D15").Value = "7.345"
$ws.Range("E15").Value = "  +3.27%  "
$ws.Range("D16").Value = "1.790.29"
$ws.Range("E16").Value = "  +0.27%  "
$ws.Range("D17").Value = "0.00001077"
$ws.Range("E17").Value = "  -1.38%  "
$ws.Range("D18").Value = "0.06677"
$ws.Range("E18").Value = "  -0.15%  "
$ws.Range("D19").Value = "82.38"
$ws.Range("D20").Value = "1.003"
$ws.Range("E20").Value = "  +0.28%  "
$ws.Range("D21").Value = "17.41"
$ws.Range("E21").Value = "  -1.49%  "
$ws.Range("D22").Value = "6.414"
$ws.Range("E22").Value = "  -1.24%  "
$ws.Range("D23").Value = "28.180.01"
$ws.Range("E23").Value = "  +3.29%  "
$ws.Range("D24").Value = "12.08"
$ws.Range("E24").Value = "  -2.39%  "
$ws.Range("D25").Value = "2.375"
$ws.Range("E25").Value = "  -0.27%  "
$ws.Range("D26").Value = "20.72"
$ws.Range("E26").Value = "  -1.97%  "
$ws.Range("D27").Value = "1.420"
$ws.Range("E27").Value = "  -3.29%  "
$ws.Range("D30").Value = "1.989.30"
$ws.Range("E30").Value = "  +0.07%  "
$ws.Range("D31").Value = "134.81"
$ws.Range("E31").Value = "  -0.83%  "
$ws.Range("D32").Value = "4.023"
$ws.Range("E32").Value = "  +0.06%  "
$ws.Range("D33").Value = "6.091"
$ws.Range("E33").Value = "  +2.36%  "
$ws.Range("D34").Value = "0.08942"
$ws.Range("E34").Value = "  +1.04%  "
$ws.Range("D35").Value = "12.74"
$ws.Range("E35").Value = "  -1.71%  "
$ws.Range("E36").Value = "  -1.07%  "
$ws.Range("D37").Value = "0.6844"
$ws.Range("E37").Value = "  +0.42%  "
$ws.Range("D38").Value = "5.347"
$ws.Range("E38").Value = "  -1.07%  "
$ws.Range("D39").Value = "0.06377"
$ws.Range("E39").Value = "  -0.93%  "
$ws.Range("D40").Value = "0.2158"
$ws.Range("E40").Value = "  -1.89%  "
$ws.Range("D41").Value = "1.248"
$ws.Range("E41").Value = "  -0.38%  "
$ws.Range("D42").Value = "1.502"
$ws.Range("E42").Value = "  -7.33%  "
$ws.Range("D43").Value = "8.304"
$ws.Range("E43").Value = "  -1.21%  "
$ws.Range("D44").Value = "14.17"
$ws.Range("E44").Value = "  -1.74%  "
$ws.Range("D45").Value = "1.003"
$ws.Range("E45").Value = "  +0.33%  "
$ws.Range("D46").Value = "0.6273"
$ws.Range("E46").Value = "  -1.80%  "
$ws.Range("D47").Value = "3.878"
$ws.Range("E47").Value = "  -0.02%  "
$ws.Range("D48").Value = "132.70"
$ws.Range("E48").Value = "  -0.03%  "
$ws.Range("D49").Value = "2.078"
$ws.Range("D50").Value = "0.07497"
$ws.Range("E50").Value = "  +5.16%  "
$ws.Range("D51").Value = "1.190"
$ws.Range("E51").Value = "  +2.32%  "

# --- Rows 28 and 29 swap (Monero <-> LidoDAOToken) ---
$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").Value = "2.408"
$ws.Range("E28").Value = "  -3.42%  "

$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").Value = "154.21"
$ws.Range("E29").Value = "  -2.04%  "
